$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal string value even when it looks like a number
# (e.g. "1.00", "0.619") without leaving a visible "text number" style on
# the cell - mirrors typing an apostrophe-prefixed value in the UI and then
# clearing the resulting cell style back to Normal.
function Set-TextValue($range, $text) {
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

$ws.Range('D2').Value = '70.191.90'
$ws.Range('E2').Value = '  +0.72%  '
$ws.Range('D3').Value = '3.559.31'
$ws.Range('E3').Value = '  +1.08%  '
Set-TextValue $ws.Range('D4') '1.00'
$ws.Range('E4').Value = '  -0.04%  '
Set-TextValue $ws.Range('D5') '604.87'
$ws.Range('E5').Value = '  +3.04%  '
Set-TextValue $ws.Range('D6') '186.01'
$ws.Range('E6').Value = '  +1.44%  '
$ws.Range('D7').Value = '3.556.63'
$ws.Range('E7').Value = '  +1.29%  '
Set-TextValue $ws.Range('D8') '0.619'
$ws.Range('E8').Value = '  +1.20%  '
$ws.Range('E9').Value = '  -0.08%  '
$ws.Range('E10').Value = '  +10.21%  '
Set-TextValue $ws.Range('D11') '0.645'
$ws.Range('E11').Value = '  +0.35%  '
Set-TextValue $ws.Range('D12') '53.87'
$ws.Range('E12').Value = '  -0.35%  '
Set-TextValue $ws.Range('D13') '0.0000309'
$ws.Range('E13').Value = '  +1.89%  '
Set-TextValue $ws.Range('D14') '9.47'
$ws.Range('E14').Value = '  +0.25%  '
$ws.Range('D15').Value = '4.125.49'
$ws.Range('E15').Value = '  +1.09%  '
$ws.Range('D16').Value = '70.305.88'
$ws.Range('E16').Value = '  +0.87%  '
$ws.Range('D17').Value = '3.561.74'
$ws.Range('E17').Value = '  +1.53%  '
Set-TextValue $ws.Range('D18') '12.66'
$ws.Range('E18').Value = '  +2.48%  '
Set-TextValue $ws.Range('D19') '18.94'
$ws.Range('E19').Value = '  -1.87%  '
Set-TextValue $ws.Range('D20') '575.48'
$ws.Range('E20').Value = '  +5.95%  '
Set-TextValue $ws.Range('D21') '0.120'
$ws.Range('E21').Value = '  +0.74%  '
Set-TextValue $ws.Range('D22') '0.993'
$ws.Range('E22').Value = '  -1.54%  '
$ws.Range('E23').Value = '  -2.31%  '
Set-TextValue $ws.Range('D24') '4.70'
$ws.Range('E24').Value = '  +3.25%  '
Set-TextValue $ws.Range('D25') '4.87'
$ws.Range('E25').Value = '  +0.82%  '
Set-TextValue $ws.Range('D26') '94.02'
$ws.Range('E26').Value = '  -1.72%  '
Set-TextValue $ws.Range('D27') '2.93'
$ws.Range('E27').Value = '  -1.35%  '
Set-TextValue $ws.Range('D28') '10.92'
$ws.Range('E28').Value = '  -1.97%  '
Set-TextValue $ws.Range('D29') '9.32'
$ws.Range('E29').Value = '  +3.00%  '
Set-TextValue $ws.Range('D30') '32.31'
$ws.Range('E30').Value = '  +0.88%  '
Set-TextValue $ws.Range('D31') '7.04'
$ws.Range('E31').Value = '  -2.73%  '
Set-TextValue $ws.Range('D32') '12.21'
$ws.Range('E32').Value = '  -1.42%  '
Set-TextValue $ws.Range('D33') '0.114'
$ws.Range('E33').Value = '  +1.97%  '
Set-TextValue $ws.Range('D34') '63.73'
$ws.Range('E34').Value = '  -0.84%  '
$ws.Range('E35').Value = '  +21.23%  '
$ws.Range('E36').Value = '  +2.60%  '
Set-TextValue $ws.Range('D37') '524.95'
$ws.Range('E37').Value = '  -3.83%  '
Set-TextValue $ws.Range('D38') '0.404'
$ws.Range('E38').Value = '  -1.34%  '
$ws.Range('B39').Value = 'Dai'
$ws.Range('C39').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue $ws.Range('D39') '0.999'
$ws.Range('E39').Value = '  +0.00%  '
$ws.Range('B40').Value = 'Maker'
$ws.Range('C40').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D40').Value = '3.627.84'
$ws.Range('E40').Value = '  +7.68%  '
Set-TextValue $ws.Range('D41') '37.43'
$ws.Range('E41').Value = '  -1.35%  '
$ws.Range('D42').Value = '0.0₃0782'
$ws.Range('E42').Value = '  +3.16%  '
Set-TextValue $ws.Range('D43') '3.51'
$ws.Range('E43').Value = '  +4.10%  '
Set-TextValue $ws.Range('D44') '0.137'
$ws.Range('E44').Value = '  +2.25%  '
Set-TextValue $ws.Range('D45') '0.0456'
$ws.Range('E45').Value = '  +4.07%  '
Set-TextValue $ws.Range('D46') '3.46'
$ws.Range('E46').Value = '  -2.54%  '
Set-TextValue $ws.Range('D47') '2.94'
$ws.Range('E47').Value = '  -0.88%  '
$ws.Range('E48').Value = '  +2.69%  '
Set-TextValue $ws.Range('D49') '9.17'
$ws.Range('E49').Value = '  +1.00%  '
$ws.Range('E50').Value = '  +0.17%  '
Set-TextValue $ws.Range('D51') '135.57'
$ws.Range('E51').Value = '  -0.83%  '
